$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "last refreshed" timestamp banner in A1 ---
$ws.Range("A1").Value = "Datos actualizados a 12 de Abril de 2020 a las 17:22"

# Row 4: Estados Unidos - updated daily figures
$ws.Range("B4").Value = 535173
$ws.Range("C4").Value = 2294
$ws.Range("D4").Value = 30604
$ws.Range("E4").Value = 483920
$ws.Range("G4").Value = 72
$ws.Range("H4").Value = 20649

# Row 15: Suiza - updated daily figures
$ws.Range("B15").Value = 25328
$ws.Range("C15").Value = 221
$ws.Range("E15").Value = 12139
$ws.Range("G15").Value = 53
$ws.Range("H15").Value = 1089

# Row 30: country re-sorted into this slot -> now "Polonia" (was "Noruega")
$ws.Range("A30").Value = "Polonia"
$ws.Range("B30").Value = 6674
$ws.Range("C30").Value = 318
$ws.Range("D30").Value = 439
$ws.Range("E30").Value = 6003
$ws.Range("F30").Value = 160
$ws.Range("G30").Value = 24
$ws.Range("H30").Value = 232

# Row 31: country re-sorted into this slot -> now "Noruega" (was "Polonia")
$ws.Range("A31").Value = "Noruega"
$ws.Range("B31").Value = 6485
$ws.Range("C31").Value = 76
$ws.Range("D31").Value = 32
$ws.Range("E31").Value = 6329
$ws.Range("F31").Value = 59
$ws.Range("G31").Value = 5
$ws.Range("H31").Value = 124

# Row 36: Pakistan - updated daily figures
$ws.Range("B36").Value = 5183
$ws.Range("C36").Value = 172
$ws.Range("E36").Value = 4069

# Row 54: Argentina - updated daily figures
$ws.Range("F54").Value = 83

# Row 55: Grecia - updated daily figures
$ws.Range("B55").Value = 2114
$ws.Range("C55").Value = 33
$ws.Range("E55").Value = 1747

# Row 75: Kazajistan - updated daily figures
$ws.Range("B75").Value = 927
$ws.Range("C75").Value = 62
$ws.Range("E75").Value = 818

# Row 82: Bulgaria - updated daily figures
$ws.Range("F82").Value = 36

# Row 85: country re-sorted into this slot -> now "Republica de Chipre" (was "Libano")
$ws.Range("A85").Value = "Republica de Chipre"
$ws.Range("B85").Value = 633
$ws.Range("C85").Value = 17
$ws.Range("D85").Value = 61
$ws.Range("E85").Value = 562
$ws.Range("F85").Value = 8
$ws.Range("H85").Value = 10

# Row 86: country re-sorted into this slot -> now "Libano" (was "Banglades")
$ws.Range("A86").Value = "Libano"
$ws.Range("B86").Value = 630
$ws.Range("C86").Value = 11
$ws.Range("D86").Value = 77
$ws.Range("E86").Value = 533
$ws.Range("F86").Value = 34
$ws.Range("G86").Value = 0
$ws.Range("H86").Value = 20

# Row 87: country re-sorted into this slot -> now "Banglades" (was "Cuba")
$ws.Range("A87").Value = "Banglades"
$ws.Range("B87").Value = 621
$ws.Range("C87").Value = 139
$ws.Range("D87").Value = 39
$ws.Range("E87").Value = 548
$ws.Range("F87").Value = 1
$ws.Range("G87").Value = 4
$ws.Range("H87").Value = 34

# Row 88: country re-sorted into this slot -> now "Cuba" (was "Republica de Chipre")
$ws.Range("A88").Value = "Cuba"
$ws.Range("B88").Value = 620
$ws.Range("D88").Value = 77
$ws.Range("E88").Value = 527
$ws.Range("F88").Value = 11
$ws.Range("H88").Value = 16

# Row 90: Oman - updated daily figures
$ws.Range("E90").Value = 486
$ws.Range("G90").Value = 1
$ws.Range("H90").Value = 4

# Row 105: Mauricio - updated daily figures
$ws.Range("B105").Value = 324
$ws.Range("C105").Value = 5
$ws.Range("D105").Value = 42
$ws.Range("E105").Value = 273
